$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns AL, AM, AN
$ws.Range("AL1").Value = "r"
$ws.Range("AM1").Value = "xo"
$ws.Range("AN1").Value = "yo"

# Copy style from an existing header cell (AK1) to the new header cells
$ws.Range("AK1").Copy()
$ws.Range("AL1:AN1").PasteSpecial(-4122) # xlPasteFormats

# Fill values for rows 2-21
$lastRow = 21
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 38).Value = 120  # AL
    $ws.Cells.Item($r, 39).Value = 145  # AM
    $ws.Cells.Item($r, 40).Value = 120  # AN
}
